$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated iteration results (rows 2-7, columns B:D only - column A is
# unchanged for those rows) plus a brand new row 8, all stored as text
# values (matching the original inline-string cell type). The leading
# apostrophe forces Excel to keep the numeric-looking text as a string
# instead of converting it to a Number.
$data = @{
    "B2" = "-0.5";                  "C2" = "3.375";                 "D2" = "16936.2850417816"
    "B3" = "-2.95206698574524e-05"; "C3" = "2.00008856200955";      "D3" = "1.00004427663443"
    "B4" = "0.666732470480711";     "C4" = "0.296186631933458";     "D4" = "0.317399175173803"
    "B5" = "0.976753098196847";     "C5" = "0.0016086922758489";    "D5" = "0.0231586063769665"
    "B6" = "0.999909611297429";     "C6" = "2.45096138939971e-08";  "D6" = "9.19113512644414e-05"
    "B7" = "1.00000152278866";      "C7" = "6.95643542769631e-12";  "D7" = "7.99779189054593e-07"
    "A8" = "7";  "B8" = "1.00000072300889";      "C8" = "1.56852308919042e-12";  "D8" = "nan"
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = "'" + $data[$addr]
}
